# Updated the packet support matrix
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (spawn global entity): add "Side effects" remark ---
$ws.Range("M4").Value = "Ignored"

# --- Row 5 (spawn mob): shorten the comment (drop "no velocity") ---
$ws.Range("G5").Value = "no uuid, no metadata, no validation"

# --- Row 6 (spawn painting): mark Partial / None / None / None + comment ---
$ws.Range("B5").Copy() | Out-Null
$ws.Range("B6").PasteSpecial(-4122) | Out-Null
$ws.Range("D5:F5").Copy() | Out-Null
$ws.Range("D6:F6").PasteSpecial(-4122) | Out-Null
$ws.Range("B6").Value = "Partial"
$ws.Range("D6").Value = "None"
$ws.Range("E6").Value = "None"
$ws.Range("F6").Value = "None"
$ws.Range("G6").Value = "no metadata, no validation, creates iplayer"

# --- Row 33 (keep alive): change "Full (no validation)" to "Full" ---
$ws.Range("B13").Copy() | Out-Null
$ws.Range("B33").PasteSpecial(-4122) | Out-Null
$ws.Range("B33").Value = "Full"

# --- Rows 40,41,42 (entity relative move / look and move / look): mark Full/None/None/None ---
foreach ($r in 40,41,42) {
    $ws.Range("B13").Copy() | Out-Null
    $ws.Range("B$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("D13:F13").Copy() | Out-Null
    $ws.Range("D$r" + ":F$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("B$r").Value = "Full"
    $ws.Range("D$r").Value = "None"
    $ws.Range("E$r").Value = "None"
    $ws.Range("F$r").Value = "None"
}

# --- Row 67 (update health): mark Partial/None/None/None + comment ---
$ws.Range("B5").Copy() | Out-Null
$ws.Range("B67").PasteSpecial(-4122) | Out-Null
$ws.Range("D5:F5").Copy() | Out-Null
$ws.Range("D67:F67").PasteSpecial(-4122) | Out-Null
$ws.Range("B67").Value = "Partial"
$ws.Range("D67").Value = "None"
$ws.Range("E67").Value = "None"
$ws.Range("F67").Value = "None"
$ws.Range("G67").Value = "no food, no saturation"

# --- Row 78 (entity teleport): mark Full/None/None/None ---
$ws.Range("B13").Copy() | Out-Null
$ws.Range("B78").PasteSpecial(-4122) | Out-Null
$ws.Range("D13:F13").Copy() | Out-Null
$ws.Range("D78:F78").PasteSpecial(-4122) | Out-Null
$ws.Range("B78").Value = "Full"
$ws.Range("D78").Value = "None"
$ws.Range("E78").Value = "None"
$ws.Range("F78").Value = "None"

# --- Rows 82-85 (login packets): drop "Disconnect (login)" entry, shift the rest up ---
$a82 = $ws.Range("A82").Value2
$a83 = $ws.Range("A83").Value2
$a84 = $ws.Range("A84").Value2
$a85 = $ws.Range("A85").Value2

# Remove the "Disconnect (login)" row entirely (A & B); leave J82 untouched.
$ws.Range("A82:B82").Clear() | Out-Null

# Shift the packet-name labels up by one row.
$ws.Range("A83").Value = $a82
$ws.Range("A84").Value = $a83
$ws.Range("A85").Value = $a84
$ws.Range("A86").Value = $a85

# "Encryption request (login)" (now row 84) has no implementation marking any more.
$ws.Range("B84:F84").Clear() | Out-Null
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B84").PasteSpecial(-4122) | Out-Null

# "Login success (login)" (now row 85) becomes Full / None / None / "Ignored".
$ws.Range("B13").Copy() | Out-Null
$ws.Range("B85").PasteSpecial(-4122) | Out-Null
$ws.Range("D13:F13").Copy() | Out-Null
$ws.Range("D85:F85").PasteSpecial(-4122) | Out-Null
$ws.Range("B85").Value = "Full"
$ws.Range("D85").Value = "None"
$ws.Range("E85").Value = "None"
$ws.Range("F85").Clear() | Out-Null
$ws.Range("F85").Value = "Ignored"

# "set compression (login)" (now row 86) keeps the plain empty marking.
$ws.Range("B86").Interior.Color = 192

$ws.Range("G26").Select() | Out-Null
